$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MORT")
$ws.Range("A4").Value = "flag_negative_pt"
$ws.Range("B4").Value = "Respondent reported wrong dates leading to negative person time calculation"
$ws.Range("C4").Value = "Please follow up with Enumerator"
$ws.Range("B10").Select() | Out-Null
